# Edit the "hashtag" sheet's A1 cell: "#hashtag1" -> "#hashtag"
# (rest of the string "#hashtag2 #hashtag3 @kişi1 @kişi2" is unchanged),
# then make that sheet the active/visible tab (mirrors the author
# switching to the "hashtag" sheet to make the edit).

$wb = $excel.ActiveWorkbook

$wsHashtag = $wb.Worksheets.Item("hashtag")
$wsHashtag.Range("A1").Value = "#hashtag #hashtag2 #hashtag3 @kişi1 @kişi2"

$wsHashtag.Activate()
